$d = $word.ActiveDocument

$replacements = @(
    @("2025-08-07 Thursday", "2025-08-08 Friday"),
    @("715×4=2860", "950×2=1900"),
    @("104×8=832", "387×5=1935"),
    @("571×5=2855", "115×4=460"),
    @("878×9=7902", "555×8=4440"),
    @("370×2=740", "887×9=7983"),
    @("112×2=224", "466×9=4194"),
    @("603×9=5427", "669×8=5352"),
    @("666×5=3330", "261×9=2349"),
    @("945×7=6615", "356×3=1068"),
    @("674×6=4044", "753×2=1506"),
    @("518×6=3108", "236×4=944"),
    @("934×3=2802", "356×4=1424"),
    @("565×6=3390", "910×2=1820"),
    @("916×2=1832", "302×9=2718"),
    @("500×7=3500", "318×5=1590"),
    @("476×3=1428", "853×8=6824"),
    @("860×5=4300", "580×3=1740"),
    @("612×8=4896", "963×4=3852"),
    @("961×4=3844", "731×9=6579"),
    @("183×5=915", "779×6=4674"),
    @("711×3=2133", "416×8=3328"),
    @("349×9=3141", "169×3=507"),
    @("172×2=344", "982×2=1964"),
    @("143×8=1144", "554×7=3878"),
    @("335×9=3015", "402×4=1608")
)

foreach ($pair in $replacements) {
    $old = $pair[0]
    $new = $pair[1]
    $range = $d.Content
    $range.Find.Execute($old, $true, $false, $false, $false, $false, $true, 1, $false, $new, 2)
}
